# Update cryptocurrency price/volume figures per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.Value = '''51.023.57'
$c.Style = "Normal"

$c = $ws.Range('E2')
$c.Value = '''  -1.90%  '
$c.Style = "Normal"

$c = $ws.Range('E3')
$c.Value = '''  -2.29%  '
$c.Style = "Normal"

$c = $ws.Range('D4')
$c.Value = '''0.999'
$c.Style = "Normal"

$c = $ws.Range('E4')
$c.Value = '''  -0.16%  '
$c.Style = "Normal"

$c = $ws.Range('D5')
$c.Value = '''374.35'
$c.Style = "Normal"

$c = $ws.Range('E5')
$c.Value = '''  +5.60%  '
$c.Style = "Normal"

$c = $ws.Range('D6')
$c.Value = '''101.59'
$c.Style = "Normal"

$c = $ws.Range('E6')
$c.Value = '''  -5.47%  '
$c.Style = "Normal"

$c = $ws.Range('D7')
$c.Value = '''0.543'
$c.Style = "Normal"

$c = $ws.Range('E7')
$c.Value = '''  -3.38%  '
$c.Style = "Normal"

$c = $ws.Range('E8')
$c.Value = '''  -0.13%  '
$c.Style = "Normal"

$c = $ws.Range('E9')
$c.Value = '''  -4.97%  '
$c.Style = "Normal"

$c = $ws.Range('D10')
$c.Value = '''36.78'
$c.Style = "Normal"

$c = $ws.Range('E10')
$c.Value = '''  -3.79%  '
$c.Style = "Normal"

$c = $ws.Range('E11')
$c.Value = '''  +0.37%  '
$c.Style = "Normal"

$c = $ws.Range('E12')
$c.Value = '''  -2.76%  '
$c.Style = "Normal"

$c = $ws.Range('D13')
$c.Value = '''18.24'
$c.Style = "Normal"

$c = $ws.Range('E13')
$c.Value = '''  -5.26%  '
$c.Style = "Normal"

$c = $ws.Range('D14')
$c.Value = '''3.369.28'
$c.Style = "Normal"

$c = $ws.Range('E14')
$c.Value = '''  -2.37%  '
$c.Style = "Normal"

$c = $ws.Range('D15')
$c.Value = '''7.36'
$c.Style = "Normal"

$c = $ws.Range('E15')
$c.Value = '''  -3.50%  '
$c.Style = "Normal"

$c = $ws.Range('D16')
$c.Value = '''2.907.03'
$c.Style = "Normal"

$c = $ws.Range('E16')
$c.Value = '''  -2.44%  '
$c.Style = "Normal"

$c = $ws.Range('D17')
$c.Value = '''0.922'
$c.Style = "Normal"

$c = $ws.Range('E17')
$c.Value = '''  -7.60%  '
$c.Style = "Normal"

$c = $ws.Range('D18')
$c.Value = '''50.981.64'
$c.Style = "Normal"

$c = $ws.Range('E18')
$c.Value = '''  -2.11%  '
$c.Style = "Normal"

$c = $ws.Range('E19')
$c.Value = '''  -6.79%  '
$c.Style = "Normal"

$c = $ws.Range('E20')
$c.Value = '''  -3.84%  '
$c.Style = "Normal"

$c = $ws.Range('D21')
$c.Value = '''12.84'
$c.Style = "Normal"

$c = $ws.Range('E21')
$c.Value = '''  -5.71%  '
$c.Style = "Normal"

$c = $ws.Range('D22')
$c.Value = '''0.0₃0941'
$c.Style = "Normal"

$c = $ws.Range('E22')
$c.Value = '''  -3.32%  '
$c.Style = "Normal"

$c = $ws.Range('D24')
$c.Value = '''259.45'
$c.Style = "Normal"

$c = $ws.Range('E24')
$c.Value = '''  -1.47%  '
$c.Style = "Normal"

$c = $ws.Range('E25')
$c.Value = '''  -1.41%  '
$c.Style = "Normal"

$c = $ws.Range('E26')
$c.Value = '''  -6.04%  '
$c.Style = "Normal"

$c = $ws.Range('E27')
$c.Value = '''  +0.03%  '
$c.Style = "Normal"

$c = $ws.Range('D28')
$c.Value = '''4.10'
$c.Style = "Normal"

$c = $ws.Range('E28')
$c.Value = '''  -4.55%  '
$c.Style = "Normal"

$c = $ws.Range('D29')
$c.Value = '''25.63'
$c.Style = "Normal"

$c = $ws.Range('E29')
$c.Value = '''  -4.62%  '
$c.Style = "Normal"

$c = $ws.Range('D30')
$c.Value = '''7.07'
$c.Style = "Normal"

$c = $ws.Range('E30')
$c.Value = '''  -6.39%  '
$c.Style = "Normal"

$c = $ws.Range('E31')
$c.Value = '''  -7.07%  '
$c.Style = "Normal"

$c = $ws.Range('D32')
$c.Value = '''6.31'
$c.Style = "Normal"

$c = $ws.Range('E32')
$c.Value = '''  +3.84%  '
$c.Style = "Normal"

$c = $ws.Range('D33')
$c.Value = '''9.83'
$c.Style = "Normal"

$c = $ws.Range('E33')
$c.Value = '''  -4.50%  '
$c.Style = "Normal"

$c = $ws.Range('E34')
$c.Value = '''  -3.77%  '
$c.Style = "Normal"

$c = $ws.Range('E35')
$c.Value = '''  +1.05%  '
$c.Style = "Normal"

$c = $ws.Range('D36')
$c.Value = '''34.06'
$c.Style = "Normal"

$c = $ws.Range('E36')
$c.Value = '''  -5.91%  '
$c.Style = "Normal"

$c = $ws.Range('E37')
$c.Value = '''  +0.44%  '
$c.Style = "Normal"

$c = $ws.Range('D38')
$c.Value = '''0.0423'
$c.Style = "Normal"

$c = $ws.Range('E38')
$c.Value = '''  -5.40%  '
$c.Style = "Normal"

$c = $ws.Range('E39')
$c.Value = '''  -6.43%  '
$c.Style = "Normal"

$c = $ws.Range('D40')
$c.Value = '''16.95'
$c.Style = "Normal"

$c = $ws.Range('E40')
$c.Value = '''  -5.24%  '
$c.Style = "Normal"

$c = $ws.Range('E41')
$c.Value = '''  -4.21%  '
$c.Style = "Normal"

$c = $ws.Range('D42')
$c.Value = '''1.83'
$c.Style = "Normal"

$c = $ws.Range('E42')
$c.Value = '''  -6.86%  '
$c.Style = "Normal"

$c = $ws.Range('D43')
$c.Value = '''0.113'
$c.Style = "Normal"

$c = $ws.Range('E43')
$c.Value = '''  -3.90%  '
$c.Style = "Normal"

$c = $ws.Range('D44')
$c.Value = '''119.44'
$c.Style = "Normal"

$c = $ws.Range('E44')
$c.Value = '''  -2.06%  '
$c.Style = "Normal"

$c = $ws.Range('D45')
$c.Value = '''21.82'
$c.Style = "Normal"

$c = $ws.Range('E45')
$c.Value = '''  -3.53%  '
$c.Style = "Normal"

$c = $ws.Range('E46')
$c.Value = '''  -2.05%  '
$c.Style = "Normal"

$c = $ws.Range('D47')
$c.Value = '''2.013.01'
$c.Style = "Normal"

$c = $ws.Range('E47')
$c.Value = '''  -5.18%  '
$c.Style = "Normal"

$c = $ws.Range('E48')
$c.Value = '''  -2.45%  '
$c.Style = "Normal"

$c = $ws.Range('E49')
$c.Value = '''  -6.58%  '
$c.Style = "Normal"

$c = $ws.Range('D50')
$c.Value = '''3.200.01'
$c.Style = "Normal"

$c = $ws.Range('E50')
$c.Value = '''  -2.22%  '
$c.Style = "Normal"

$c = $ws.Range('E51')
$c.Value = '''  -2.14%  '
$c.Style = "Normal"
